$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 8-9: fill in session 7/8 expenditure breakdown ---
$ws.Range("F8").Value = 10
$ws.Range("G8").Value = 240
$ws.Range("C8").Formula = "=SUM(E8:G8)"

$ws.Range("E9").Value = 356
$ws.Range("F9").Value = 25
$ws.Range("C9").Formula = "=SUM(E9:G9)"

# --- Row 12-13: fill in session 9/10 expenditure breakdown ---
$ws.Range("E12").Value = 332
$ws.Range("F12").Value = 10
$ws.Range("G12").Formula = "= 8+4+2+2+4.5+8.5+7+5+14+10+9+6+10+5+3+4+5.5+4.5+11+2+3+7+2"
$ws.Range("C12").Formula = "=SUM(E12:G12)"

$ws.Range("E13").Value = 309.5
$ws.Range("F13").Value = 30
$ws.Range("G13").Formula = "=59 + 53 + 41.5"
$ws.Range("C13").Formula = "=SUM(E13:G13)"

# --- C16: "Yes." answer to "More cash?" ---
$ws.Range("C16").Value = "Yes."

# --- New payoff calculation table (rows 29-36) ---
$ws.Range("A29").Value = "Friday, end of the week"
$ws.Range("E29").Value = "Ordered"
$ws.Range("G29").Value = "Actual ( correction)"
$ws.Range("H29").Value = "Needed:"

$ws.Range("A30").Value = 20
$ws.Range("B30").Value = 26
$ws.Range("C30").Formula = "=A30*B30"
$ws.Range("E30").Value = "\item 50 * 20 euro = 1 000 euro"

$ws.Range("A31").Value = 10
$ws.Range("B31").Value = 58
$ws.Range("C31").Formula = "=A31*B31"
$ws.Range("E31").Value = "\item 160 * 10 euro = 1 600 euro"
$ws.Range("H31").Value = "50 kertaa kymmenen euroa"
$ws.Range("J31").Formula = "=50*10"

$ws.Range("A32").Value = 5
$ws.Range("B32").Value = 65
$ws.Range("C32").Formula = "=A32*B32"
$ws.Range("E32").Value = "\item 160 * 5 euro = 800 euro"
$ws.Range("H32").Value = "50 kertaa viisi euroa"
$ws.Range("J32").Formula = "=50*5"

$ws.Range("A33").Value = 2
$ws.Range("B33").Value = 63
$ws.Range("C33").Formula = "=A33*B33"
$ws.Range("E33").Value = "\item 200 * 2 euro = 400 euro"
$ws.Range("H33").Value = "100 kertaa 2 euroa"

$ws.Range("A34").Value = 1
$ws.Range("B34").Value = 197
$ws.Range("C34").Formula = "=A34*B34"
$ws.Range("E34").Value = "\item 180 * 1 euro = 180 euro"
$ws.Range("G34").Value = "250 * 1"

$ws.Range("A35").Value = 36
$ws.Range("B35").Value = 0.5
$ws.Range("C35").Formula = "=A35*B35"
$ws.Range("E35").Value = "\item 40 * 50 snt = 20 euro"

$ws.Range("C36").Formula = "=SUM(C30:C35)"

$ws.Range("E30").Select()

$wb.Save()
